# Apply the change: insert a new data row (row 156) into the worksheet,
# shifting the existing rows 156-226 down to 157-227, then populate the
# new row 156 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156; this shifts rows 156:226 down to 157:227
# and extends the sheet dimension to A1:R227 automatically.
$ws.Rows(156).Insert()

# Populate the newly inserted row 156 with its data values.
$ws.Range("A156").Value = 10
$ws.Range("B156").Value = "Vega Modelo de Temuco"
$ws.Range("C156").Value = "La Araucanía"
$ws.Range("D156").Value = 44489
$ws.Range("E156").Value = 9
$ws.Range("F156").Value = 100112037
$ws.Range("G156").Value = "Cebollín"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 160
$ws.Range("K156").Value = 5000
$ws.Range("L156").Value = 5000
$ws.Range("M156").Value = 5000
$ws.Range("N156").Value = "`$/docena de paquetes"
$ws.Range("O156").Value = "Región de O'Higgins"
$ws.Range("P156").Value = 417
$ws.Range("Q156").Value = 12
$ws.Range("R156").Value = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the
# other rows in column D.
$ws.Range("D156").NumberFormat = $ws.Range("D157").NumberFormat
